# Add a new column W ("position_idx") to Sheet1 holding the 0-based row
# position of each data row (row 2 -> 0, row 3 -> 1, ... row 260 -> 258),
# matching the PSM "matching without replacement" fix described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell W1, formatted like the other header cells (bold, border,
#     centered/top-aligned) by copying the format from V1 -----------------
$ws.Range("W1").Value = "position_idx"
$ws.Range("V1").Copy() | Out-Null
$ws.Range("W1").PasteSpecial(-4122) | Out-Null

# --- Data cells W2:W260 = 0,1,2,...,258 ----------------------------------
$lastRow = 260
$firstRow = 2
$n = $lastRow - $firstRow + 1

$values = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $values[$i,0] = $i
}

$ws.Range("W2:W260").Value = $values
